# Update "Bibi" annual billing data for the 2025 row (row 9) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 corresponds to Ano = 2025.
$ws.Range("B9").Value = 3655917.35    # Faturamento em Produtos
$ws.Range("C9").Value = 575985.24     # Faturamento em Servicos
$ws.Range("D9").Value = 4231902.59    # Total

$ws.Range("E9").Value = 13.61055052072926    # Faturamento em Servicos (%)
$ws.Range("F9").Value = 86.38944947927074    # Faturamento em Produtos (%)

$ws.Range("G9").Value = -44.33389854881037   # Evolucao Faturamento em Servicos (%)
$ws.Range("H9").Value = -33.97918601278863   # Evolucao Faturamento em Produtos (%)

$ws.Range("I9").Value = 36973   # Qtd Produtos
$ws.Range("J9").Value = 1575    # Qtd Servicos
$ws.Range("K9").Value = 38548   # Total Itens
$ws.Range("L9").Value = 26628   # Qtd Vendas

$ws.Range("M9").Value = 158.9267909719093    # Ticket Medio Anual
$ws.Range("N9").Value = 8.502557174396941    # Evolucao Ticket Medio (%)
